# "Journal de travail" update
# - Fills in the last logged entry (2018-02-16) with its duration + description.
# - Adds a new log entry right after it (same day) for documentation work.
# - The "Log" table (and its AutoFilter) grows by one row to keep wrapping the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Agenda")

# Insert a brand-new row at 20 -> everything that used to be on/after row 20
# (rows 20-41) shifts down to rows 21-42.
$ws.Rows.Item(20).Insert()
$ws.Rows.Item(20).RowHeight = 30.75

# Grow the "Log" table so the new row (and the row below it) stay part of it.
$lo = $ws.ListObjects.Item("Log")
$lo.Resize($ws.Range("B2:E42"))

# Row 19 (2018-02-16, already had a date) now gets its duration + description.
$ws.Range("C19").Value = "45 min"
$ws.Range("D19").Value = "Amélioration + finitions de l'interface graphique de la génération de terrain"

# Row 20 is the new entry: same date, documentation work.
$ws.Range("B20").Value = 43147
$ws.Range("C20").Value = "45 min"
$ws.Range("D20").Value = "Rédaction de la documentation"

# Reflect the new scroll position / selection used while editing.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("B21").Select()
